# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from serial date 45186 to serial date 45188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45188
}
